$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demands")
$ws.Columns("F:F").Insert()
